# "cambio de fracciones e historico" -- update the reporting quarter dates
# (and the matching "fecha de validacion"/"fecha de actualizacion" dates)
# for the three data rows of the "Reporte de Formatos" sheet, and move the
# active selection the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Row 8 - Inscripcion UPP
$ws.Range("B8").Value  = 44743   # periodo que se informa (inicio) -> 01/07/2022
$ws.Range("C8").Value  = 44834   # periodo que se informa (fin)    -> 30/09/2022
$ws.Range("AA8").Value = 44844   # fecha de validacion             -> 10/10/2022
$ws.Range("AB8").Value = 44844   # fecha de actualizacion          -> 10/10/2022

# Row 9 - Becas Institucionales UPP
$ws.Range("B9").Value  = 44743
$ws.Range("C9").Value  = 44834
$ws.Range("AA9").Value = 44844
$ws.Range("AB9").Value = 44844

# Row 10 - Estadia Profesional UPP
$ws.Range("B10").Value  = 44743
$ws.Range("C10").Value  = 44834
$ws.Range("AA10").Value = 44844
$ws.Range("AB10").Value = 44844

# Leave the sheet scrolled/selected the way it was saved: viewport starting
# around row 6 with C8 as the active cell.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select() | Out-Null
